$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 537
$ws.Range("B2").Value = 377
$ws.Range("B3").Value = 687
$ws.Range("B4").Value = 468
$ws.Range("B5").Value = 641
$ws.Range("B6").Value = 436
$ws.Range("B7").Value = 572
